# Update countries & provincias Spain
# Refresh the "Pais" COVID data table: re-sort/reshuffle a handful of
# country rows (new case numbers caused them to swap ranking order) and
# bump the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Timestamp footer (row 1 / A1) ----
$ws.Range("A1").Value = "Datos actualizados a 24 de Julio de 2020 a las 04:46"

# ---- Country data rows that changed (name and/or totals) ----
# Each row: Country(A), Casos totales(B), Nuevos casos(C), Casos activos(D),
#           Recuperados(E), Casos criticos(F), Muertes hoy(G), Muertes(H)

# Belgica / Bolivia swap ranking (rows 36-37)
$ws.Range("A36").Value = "Bolivia"
$ws.Range("B36").Value = 65252
$ws.Range("C36").Value = 1117
$ws.Range("D36").Value = 20030
$ws.Range("E36").Value = 42815
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 79
$ws.Range("H36").Value = 2407

$ws.Range("A37").Value = "Belgica"
$ws.Range("B37").Value = 64627
$ws.Range("C37").Value = 0
$ws.Range("D37").Value = 17347
$ws.Range("E37").Value = 37472
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 9808

# Corea del Sur data refresh (row 72, name unchanged)
$ws.Range("A72").Value = "Corea del Sur"
$ws.Range("B72").Value = 13979
$ws.Range("C72").Value = 41
$ws.Range("D72").Value = 12817
$ws.Range("E72").Value = 864
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 298

# Dinamarca / El Salvador / Australia 3-way re-rank (rows 74-76)
$ws.Range("A74").Value = "Australia"
$ws.Range("B74").Value = 13595
$ws.Range("C74").Value = 289
$ws.Range("D74").Value = 8775
$ws.Range("E74").Value = 4681
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 6
$ws.Range("H74").Value = 139

$ws.Range("A75").Value = "Dinamarca"
$ws.Range("B75").Value = 13390
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 12299
$ws.Range("E75").Value = 479
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 612

$ws.Range("A76").Value = "El Salvador"
$ws.Range("B76").Value = 13377
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 7276
$ws.Range("E76").Value = 5729
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 372

# Camboya data refresh (row 175, name unchanged)
$ws.Range("A175").Value = "Camboya"
$ws.Range("B175").Value = 202
$ws.Range("C175").Value = 4
$ws.Range("D175").Value = 142
$ws.Range("E175").Value = 60
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

# Gibraltar data refresh (row 178, name unchanged)
$ws.Range("A178").Value = "Gibraltar"
$ws.Range("B178").Value = 184
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 180
$ws.Range("E178").Value = 4
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# Aruba data refresh (row 183, name unchanged)
$ws.Range("A183").Value = "Aruba"
$ws.Range("B183").Value = 117
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 101
$ws.Range("E183").Value = 13
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 3

# Liechtenstein / Islas Turcas y Caicos swap ranking (rows 188-189)
$ws.Range("A188").Value = "Islas Turcas y Caicos"
$ws.Range("B188").Value = 90
$ws.Range("C188").Value = 4
$ws.Range("D188").Value = 26
$ws.Range("E188").Value = 62
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 2

$ws.Range("A189").Value = "Liechtenstein"
$ws.Range("B189").Value = 86
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 81
$ws.Range("E189").Value = 4
$ws.Range("F189").Value = 0
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 1

# San Martin (Parte Francesa) / Macao / Belice 3-way re-rank (rows 194-196)
$ws.Range("A194").Value = "Belice"
$ws.Range("B194").Value = 47
$ws.Range("C194").Value = 4
$ws.Range("D194").Value = 25
$ws.Range("E194").Value = 20
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 2

$ws.Range("A195").Value = "San Martin (Parte Francesa)"
$ws.Range("B195").Value = 46
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 39
$ws.Range("E195").Value = 4
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 3

$ws.Range("A196").Value = "Macao"
$ws.Range("B196").Value = 46
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 46
$ws.Range("E196").Value = 0
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

# Papua Nueva Guinea data refresh (row 199, name unchanged)
$ws.Range("A199").Value = "Papua Nueva Guinea"
$ws.Range("B199").Value = 31
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 11
$ws.Range("E199").Value = 20
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0

# Nueva Caledonia data refresh (row 205, name unchanged)
$ws.Range("A205").Value = "Nueva Caledonia"
$ws.Range("B205").Value = 22
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 22
$ws.Range("E205").Value = 0
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

# Islas Malvinas / Groenlandia swap ranking (rows 210-211).
# Their totals (B-H) are identical for both rows already, so only the
# country names need to change places.
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"
